$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 310 (pushing the
# existing rows 310-408 down to 311-409).
$ws.Rows(310).Insert()

$ws.Range("A310").Value = 4
$ws.Range("B310").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C310").Value = "Los Lagos"
$ws.Range("D310").Value = 44985
$ws.Range("E310").Value = 10
$ws.Range("F310").Value = 100112040
$ws.Range("G310").Value = "Cilantro"
$ws.Range("H310").Value = "Sin especificar"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 160
$ws.Range("K310").Value = 7000
$ws.Range("L310").Value = 7000
$ws.Range("M310").Value = 7000
$ws.Range("N310").Value = "$/docena de atados (2 kilos)"
$ws.Range("O310").Value = "Región de La Araucanía"
$ws.Range("P310").Value = 3500
$ws.Range("Q310").Value = 2
$ws.Range("R310").Value = "Hortaliza"
